$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Commit: "commit model 3 output for analysis"
# Zero out the tardiness-at-customer sample values in A2:A9
# (A1 and A10 already hold 0 and remain unchanged).
$ws.Range("A2:A9").Value = 0
